$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Fill in the new time-log entries for rows 20-22 and 24-27.
# The shared-string table must end up with the new descriptions in the same
# order as the target workbook, so the D (Description) column values are
# written first, in the exact order required, before anything else touches
# those cells.
# ---------------------------------------------------------------------------

$ws.Range("D24").Value = "Setup dual boot Windows-Ubuntu at home"
$ws.Range("D25").Value = "Meeting with professor gross, cloning the hard drive"
$ws.Range("D26").Value = "Looking through config files for NIS, more research"
$ws.Range("D21").Value = "Meeting professor Gross, installing bullseye 11.7 from scratch"
$ws.Range("D20").Value = "Trying to fix broken package dependencies"
$ws.Range("D22").Value = "More research, installing ubuntu on personal laptop for SSH"
$ws.Range("D27").Value = "Demo 1 Video, Installing software to make a video"

# Row 20 - 9/19/2023, 2 hours
$ws.Range("B20").Value = 45188
$ws.Range("B20").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("B20").HorizontalAlignment = -4108
$ws.Range("C20").Value = 2

# Row 21 - 9/21/2023, 3 hours
$ws.Range("B21").Value = 45190
$ws.Range("B21").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("B21").HorizontalAlignment = -4108
$ws.Range("C21").Value = 3

# Row 22 - 9/23/2023, 4 hours
$ws.Range("B22").Value = 45192
$ws.Range("B22").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("B22").HorizontalAlignment = -4108
$ws.Range("C22").Value = 4
$ws.Range("C22").HorizontalAlignment = -4108

# Row 23 is fully cleared out (no week-5 row divider keeps it empty)
[void]$ws.Range("B23:C23").Clear()

# Row 24 - Week 5, 9/26/2023, 6 hours
$ws.Range("A24").Value = 5
$ws.Range("A24").HorizontalAlignment = -4108
$ws.Range("B24").Value = 45195
$ws.Range("B24").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("B24").HorizontalAlignment = -4108
$ws.Range("C24").Value = 6

# Row 25 - 9/28/2023, 2.5 hours
$ws.Range("B25").Value = 45197
$ws.Range("B25").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("B25").HorizontalAlignment = -4108
$ws.Range("C25").Value = 2.5

# Row 26 - 9/29/2023, 3.5 hours
$ws.Range("B26").Value = 45198
$ws.Range("B26").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("B26").HorizontalAlignment = -4108
$ws.Range("C26").Value = 3.5

# Row 27 - 10/1/2023, 3 hours
$ws.Range("B27").Value = 45200
$ws.Range("B27").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("B27").HorizontalAlignment = -4108
$ws.Range("C27").Value = 3

# ---------------------------------------------------------------------------
# Misc. workbook/view state updates
# ---------------------------------------------------------------------------

# Move the active selection to D27 (last entry added)
[void]$ws.Range("D27").Select()

# Set the print orientation to portrait
$ws.PageSetup.Orientation = 1

Write-Output "done"
